# Updated cryptos list on Sun Feb 25 04:28:56 UTC 2024 with GitHub Actions
# Applies refreshed Price (D) / Volume(1h) (E) values, and fixes the row order
# for a few coins (VeChain/OKB/Toncoin and Stacks/ARBITRUM).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.555.71'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.016.27'
$ws.Range('E3').Value = '  +2.71%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''378.29'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '''103.19'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').Value = '''0.545'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.595'
$ws.Range('E9').Value = '  +2.75%  '
$ws.Range('D10').Value = '''36.70'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '''0.0860'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').Value = '3.500.08'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').Value = '''18.51'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('D15').Value = '''7.74'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '3.016.25'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').Value = '''0.981'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '''10.49'
$ws.Range('E18').Value = '  -14.27%  '
$ws.Range('D19').Value = '51.567.06'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').Value = '''3.03'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = '''12.48'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '''69.95'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').Value = '''267.39'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').Value = '''3.13'
$ws.Range('E25').Value = '  -3.31%  '
$ws.Range('D26').Value = '''8.21'
$ws.Range('E26').Value = '  +3.39%  '
$ws.Range('D27').Value = '''7.53'
$ws.Range('E27').Value = '  +6.04%  '
$ws.Range('D28').Value = '''0.172'
$ws.Range('E28').Value = '  +5.82%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '''26.21'
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').Value = '''10.29'
$ws.Range('E32').Value = '  +2.90%  '
$ws.Range('D33').Value = '''34.17'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = '''0.0455'
$ws.Range('E34').Value = '  +5.50%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''50.67'
$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').Value = '''2.05'
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = '''3.29'
$ws.Range('E38').Value = '  +6.92%  '
$ws.Range('D39').Value = '''17.27'
$ws.Range('E39').Value = '  +4.54%  '
$ws.Range('E40').Value = '  +10.78%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''2.59'
$ws.Range('E41').Value = '  +4.67%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''1.86'
$ws.Range('E42').Value = '  +2.92%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').Value = '''126.59'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('E45').Value = '  +8.59%  '
$ws.Range('D46').Value = '''21.64'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('D49').Value = '2.031.27'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').Value = '3.317.39'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').Value = '''0.0320'
$ws.Range('E51').Value = '  +1.72%  '
